$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (original values are stored as text).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '68.953.98'
$ws.Range('E2').Value = '  +1.92%  '
$ws.Range('D3').Value = '2.499.66'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D5').Value = '591.35'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = '175.26'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').Value = '2.498.21'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').Value = '0.152'
$ws.Range('E10').Value = '  +6.94%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = '4.96'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '0.336'
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').Value = '2.956.61'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').Value = '25.58'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '68.910.59'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('D17').Value = '0.0000173'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '2.509.34'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').Value = '359.32'
$ws.Range('E19').Value = '  +2.27%  '
$ws.Range('D20').Value = '7.50'
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').Value = '10.84'
$ws.Range('E21').Value = '  -1.50%  '
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '69.75'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('D25').Value = '4.15'
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').Value = '8.91'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('D27').Value = '1.66'
$ws.Range('E27').Value = '  -6.74%  '
$ws.Range('D28').Value = '2.629.28'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').Value = '1.02'
$ws.Range('E29').Value = '  +2.74%  '
$ws.Range('D30').Value = '506.13'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').Value = '0.0₃0877'
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('D32').Value = '7.71'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').Value = '  -3.78%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '1.76'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '163.18'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('D38').Value = '18.57'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').Value = '18.68'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '1.30'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('D42').Value = '1.68'
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').Value = '4.73'
$ws.Range('E43').Value = '  -2.25%  '
$ws.Range('D44').Value = '0.317'
$ws.Range('E44').Value = '  -3.85%  '
$ws.Range('E45').Value = '  -5.41%  '
$ws.Range('D46').Value = '148.96'
$ws.Range('E46').Value = '  +2.46%  '
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').Value = '0.509'
$ws.Range('E48').Value = '  -1.17%  '
$ws.Range('D49').Value = '0.0733'
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('D51').Value = '0.573'
$ws.Range('E51').Value = '  -2.26%  '
